# Change the table style ("Table Design" gallery selection) used by the
# table on slide 16 from the deck's custom "Table_0" style to the built-in
# PowerPoint table style {8FCF13EC-4D5E-44F5-9F32-9A5EC6547A87}.
#
# PowerPoint exposes the currently-applied table style through
# Table.Style (read-only display) but changing it has to go through
# Table.ApplyStyle(styleId) - that's the supported write path (simply
# assigning to .Style is not honoured).

$p = $ppt.ActivePresentation

$slide = $p.Slides.Item(16)

# Find the shape that owns the table on this slide.
$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
    }
}

$table = $tableShape.Table
$table.ApplyStyle("{8FCF13EC-4D5E-44F5-9F32-9A5EC6547A87}")
